# Applies the "Updated cryptos list" data refresh described in the commit.
# For each changed cell we set the new text value. Column D holds numeric-
# looking values (e.g. "20.10", "1.00", "48.336.96") that must stay as literal
# text (matching the source inlineStr cells), so we force the cell to Text
# format before assigning - otherwise Excel would coerce them into numbers
# and silently drop significant trailing zeros / multi-dot formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '48.336.96'

$ws.Range("E2").Value = '  +2.43%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.515.85'

$ws.Range("E3").Value = '  +1.28%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.81'

$ws.Range("E5").Value = '  +0.30%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.91'

$ws.Range("E6").Value = '  +0.82%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.528'

$ws.Range("E7").Value = '  +1.24%  '

# Row 8
$ws.Range("E8").Value = '  +0.01%  '

# Row 9
$ws.Range("E9").Value = '  +1.15%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.08'

$ws.Range("E10").Value = '  +2.68%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.10'

$ws.Range("E11").Value = '  +9.57%  '

# Row 12
$ws.Range("E12").Value = '  +1.38%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.907.93'

$ws.Range("E15").Value = '  +1.24%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.518.70'

$ws.Range("E16").Value = '  +1.31%  '

# Row 17
$ws.Range("E17").Value = '  +0.54%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '48.168.37'

$ws.Range("E18").Value = '  +2.25%  '

# Row 19
$ws.Range("E19").Value = '  -1.47%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.78'

$ws.Range("E20").Value = '  +2.63%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0953'

$ws.Range("E21").Value = '  +1.24%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.74'

$ws.Range("E22").Value = '  +0.26%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.32'

$ws.Range("E23").Value = '  +2.65%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '277.45'

$ws.Range("E24").Value = '  +12.91%  '

# Row 25
$ws.Range("E25").Value = '  +0.90%  '

# Row 26
$ws.Range("E26").Value = '  +0.04%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.94'

$ws.Range("E27").Value = '  +1.07%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.40'

$ws.Range("E28").Value = '  +4.98%  '

# Row 29
$ws.Range("E29").Value = '  -1.22%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.52'

$ws.Range("E30").Value = '  +3.09%  '

# Row 31
$ws.Range("E31").Value = '  -1.22%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.22'

$ws.Range("E32").Value = '  -0.99%  '

# Row 33
$ws.Range("E33").Value = '  -3.60%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'

$ws.Range("E35").Value = '  +0.04%  '

# Row 36
$ws.Range("E36").Value = '  +0.74%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.97'

$ws.Range("E37").Value = '  +0.87%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.64'

$ws.Range("E38").Value = '  -2.38%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.97'

$ws.Range("E39").Value = '  +0.98%  '

# Row 40
$ws.Range("B40").Value = 'Stellar'

$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.112'

$ws.Range("E40").Value = '  +0.29%  '

# Row 41
$ws.Range("B41").Value = 'Monero'

$ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '122.54'

$ws.Range("E41").Value = '  +3.19%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.55'

$ws.Range("E43").Value = '  -5.74%  '

# Row 44
$ws.Range("E44").Value = '  +3.27%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.003.08'

$ws.Range("E45").Value = '  +0.47%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.16'

$ws.Range("E46").Value = '  +5.22%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.86'

$ws.Range("E47").Value = '  +3.64%  '

# Row 48
$ws.Range("E48").Value = '  -1.18%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.04'

$ws.Range("E49").Value = '  -0.77%  '

# Row 50
$ws.Range("E50").Value = '  +3.18%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '80.32'

$ws.Range("E51").Value = '  +3.70%  '
